$d = $word.ActiveDocument
$p24 = $d.Paragraphs(24)
$p25 = $d.Paragraphs(25)
$full = $d.Range($p24.Range.Start, $p25.Range.End)
Write-Host "full start=$($full.Start) end=$($full.End) text=[$($full.Text)]"
$full.Font.LanguageID = "es-419"
Write-Host "full font lang=$($full.Font.LanguageID)"
